$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Activate()

# Rename database identifiers: add "yetl_" prefix to the ad_works_lt databases.
$used = $ws.UsedRange
$used.Replace("control_ad_works_lt", "yetl_control_ad_works_lt", 2, 1, $false, $false, $false, $false)
$used.Replace("landing_ad_works_lt", "yetl_landing_ad_works_lt", 2, 1, $false, $false, $false, $false)
$used.Replace("raw_ad_works_lt", "yetl_raw_ad_works_lt", 2, 1, $false, $false, $false, $false)
$used.Replace("base_ad_works_lt", "yetl_base_ad_works_lt", 2, 1, $false, $false, $false, $false)

$ws.Range("G20").Select()
